$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for season record: Wins, Losses, Ties
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style from an existing header cell (AC1) onto the new headers
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Restore the values after the paste-special (paste format should not touch values,
# but keep this explicit/safe)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record for every data row (2 through 43):
# every team in this sheet shares the same season record: 96 wins, 66 losses, 0 ties
for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 30).Value = 96   # column AD
    $ws.Cells.Item($r, 31).Value = 66   # column AE
    $ws.Cells.Item($r, 32).Value = 0    # column AF
}
